# Final Touches before Submission of documents
# Re-apply cell alignment formatting across the Product Backlog sheet:
#   - header row (A1:F1): center both ways
#   - ID / Sprint(most) / Priority / Status / Story Points columns: center both ways
#   - Sprint cells that hold multi-line text (C7:C9): center both ways + wrap
#   - Story column cells (B): left-align, center vertically, wrap text
# Also move the active selection to G1 (matches the saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft   = -4131
$xlCenter = -4108

# --- Header row: A1:F1 -> horizontal=center, vertical=center ---
$rng = $ws.Range("A1:F1")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

# --- Plain data columns: horizontal=center, vertical=center ---
$rng = $ws.Range("A2:A12")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("C2:C4")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("C5")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("C6")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("C10:C12")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("D2:D12")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("E2:E12")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

$rng = $ws.Range("F2:F12")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter

# --- Sprint cells with two-line values: center both ways + wrap ---
$rng = $ws.Range("C7:C9")
$rng.HorizontalAlignment = $xlCenter
$rng.VerticalAlignment = $xlCenter
$rng.WrapText = $true

# --- Story column: left-align, vertical=center, keep wrap ---
$rng = $ws.Range("B2:B4")
$rng.HorizontalAlignment = $xlLeft
$rng.VerticalAlignment = $xlCenter
$rng.WrapText = $true

$rng = $ws.Range("B5")
$rng.HorizontalAlignment = $xlLeft
$rng.VerticalAlignment = $xlCenter
$rng.WrapText = $true

$rng = $ws.Range("B6:B7")
$rng.HorizontalAlignment = $xlLeft
$rng.VerticalAlignment = $xlCenter
$rng.WrapText = $true

$rng = $ws.Range("B8:B12")
$rng.HorizontalAlignment = $xlLeft
$rng.VerticalAlignment = $xlCenter
$rng.WrapText = $true

# --- Move active selection to G1 ---
$ws.Range("G1").Select()
